$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163. Writing date/time-look-alike strings as literal
# text (not auto-converted dates) by entering them as ="..."  formulas and
# then collapsing to values via Copy + PasteSpecial(xlPasteValues). This keeps
# the cell as a plain shared-string cell (t="s") without adding any new
# number-format style to styles.xml (unlike NumberFormat="@" or a quote prefix).

$ws.Range("A2").Value = 1
$ws.Range("B2").Formula = "=""2023-07-20"""
$ws.Range("C2").Formula = "=""08:38:26"""
$bc = $ws.Range("B2:C2")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D2").Value = "test"
$ws.Range("E2").Value = "test surname"
$ws.Range("F2").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G2").Value = "lucp2284"

$ws.Range("A3").Value = 2
$ws.Range("B3").Formula = "=""2023-07-20"""
$ws.Range("C3").Formula = "=""08:44:59"""
$bc = $ws.Range("B3:C3")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D3").Value = "test"
$ws.Range("E3").Value = "test surname"
$ws.Range("F3").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G3").Value = "lucp2284"

$ws.Range("A4").Value = 3
$ws.Range("B4").Formula = "=""2023-07-20"""
$ws.Range("C4").Formula = "=""08:44:59"""
$bc = $ws.Range("B4:C4")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D4").Value = "test"
$ws.Range("E4").Value = "test surname"
$ws.Range("F4").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G4").Value = "lucp2284"

$ws.Range("A5").Value = 4
$ws.Range("B5").Formula = "=""2023-07-20"""
$ws.Range("C5").Formula = "=""08:47:45"""
$bc = $ws.Range("B5:C5")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D5").Value = "test"
$ws.Range("E5").Value = "test surname"
$ws.Range("F5").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G5").Value = "lucp2284"

$ws.Range("A6").Value = 5
$ws.Range("B6").Formula = "=""2023-07-20"""
$ws.Range("C6").Formula = "=""08:51:44"""
$bc = $ws.Range("B6:C6")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D6").Value = "test"
$ws.Range("E6").Value = "test surname"
$ws.Range("F6").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G6").Value = "lucp2284"

$ws.Range("A7").Value = 6
$ws.Range("B7").Formula = "=""2023-07-20"""
$ws.Range("C7").Formula = "=""08:53:36"""
$bc = $ws.Range("B7:C7")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D7").Value = "test"
$ws.Range("E7").Value = "test surname"
$ws.Range("F7").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G7").Value = "lucp2284"

$ws.Range("A8").Value = 7
$ws.Range("B8").Formula = "=""2023-07-20"""
$ws.Range("C8").Formula = "=""08:59:31"""
$bc = $ws.Range("B8:C8")
$bc.Copy()
$bc.PasteSpecial(-4163)
$ws.Range("D8").Value = "test"
$ws.Range("E8").Value = "test surname"
$ws.Range("F8").Value = "marc.geraerts@uhasselt.be"
$ws.Range("G8").Value = "lucp2284"

$excel.CutCopyMode = 0
